$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E contain numeric-looking values that must remain text (matching
# the original inline-string cells). Force Text format before writing so COM
# doesn't silently coerce strings like '82.89' into numbers, then drop back to
# the Normal style so no stray formatting is left behind on the cells.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '43.613.99'
$ws.Range("E2").Value = '  +2.70%  '

$ws.Range("D3").Value = '2.191.11'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = '257.68'
$ws.Range("E5").Value = '  +1.30%  '

$ws.Range("D6").Value = '82.89'
$ws.Range("E6").Value = '  +11.67%  '

$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  +1.08%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").Value = '0.589'
$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").Value = '43.71'
$ws.Range("E10").Value = '  +8.30%  '

$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("D12").Value = '7.04'
$ws.Range("E12").Value = '  +3.68%  '

$ws.Range("E13").Value = '  +1.84%  '

$ws.Range("D14").Value = '2.520.56'
$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("D15").Value = '14.31'
$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("D16").Value = '2.190.71'
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").Value = '0.777'
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").Value = '43.530.16'
$ws.Range("E18").Value = '  +2.68%  '

$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.90'
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").Value = '69.39'
$ws.Range("E21").Value = '  -2.29%  '

$ws.Range("D22").Value = '2.35'
$ws.Range("E22").Value = '  +10.31%  '

$ws.Range("D23").Value = '230.42'
$ws.Range("E23").Value = '  +2.00%  '

$ws.Range("D24").Value = '8.78'
$ws.Range("E24").Value = '  -7.64%  '

$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("D26").Value = '41.98'
$ws.Range("E26").Value = '  +11.16%  '

$ws.Range("D27").Value = '10.62'
$ws.Range("E27").Value = '  +1.11%  '

$ws.Range("D28").Value = '3.38'
$ws.Range("E28").Value = '  -0.08%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +3.16%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +1.84%  '

$ws.Range("D31").Value = '173.54'
$ws.Range("E31").Value = '  +0.96%  '

$ws.Range("D32").Value = '20.34'
$ws.Range("E32").Value = '  +1.27%  '

$ws.Range("D33").Value = '0.0868'
$ws.Range("E33").Value = '  +5.39%  '

$ws.Range("D34").Value = '5.35'
$ws.Range("E34").Value = '  +3.55%  '

$ws.Range("D35").Value = '0.115'
$ws.Range("E35").Value = '  +6.00%  '

$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.0359'
$ws.Range("E37").Value = '  +5.16%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '4.45'
$ws.Range("E38").Value = '  +5.50%  '

$ws.Range("D39").Value = '12.48'
$ws.Range("E39").Value = '  +2.59%  '

$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  +9.63%  '

$ws.Range("D41").Value = '2.09'
$ws.Range("E41").Value = '  +1.11%  '

$ws.Range("D42").Value = '62.73'
$ws.Range("E42").Value = '  +5.70%  '

$ws.Range("D43").Value = '5.46'
$ws.Range("E43").Value = '  +5.78%  '

$ws.Range("D44").Value = '0.198'
$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("D45").Value = '100.39'
$ws.Range("E45").Value = '  -1.55%  '

$ws.Range("D46").Value = '0.0974'
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("E48").Value = '  +3.06%  '

$ws.Range("E49").Value = '  +1.68%  '

$ws.Range("D50").Value = '0.440'
$ws.Range("E50").Value = '  -4.26%  '

$ws.Range("D51").Value = '1.47'
$ws.Range("E51").Value = '  +17.49%  '

$priceRange.Style = "Normal"
